$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 8) mirroring the existing data rows.
$row = 8

$ws.Cells.Item($row, 1).Value = 42608.888819444444
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = 68
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"
